$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Green Line")
$ws.Rows.Item(2).Delete()
